# Road_Map.xlsx update: condensed table (col D) now depends on uniform
# skus, plus the rest of the ProFocus websites added to column A, and the
# model names in column B corrected/normalized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: normalize / correct model names -----------------------------
$ws.Range("B5").Value  = "Epson T5470"
$ws.Range("B6").Value  = "Epson T3270"
$ws.Range("B7").Value  = "Epson T5270"
$ws.Range("B8").Value  = "Epson T7270"
$ws.Range("B10").Value = "Canon IPF670e"

# --- Column D: condensed / "fixed" table ------------------------------------
$ws.Range("D2").Value  = "Adorama (fixed)"
$ws.Range("D3").Value  = "GoVets (fixed)"
$ws.Range("D4").Value  = "Plotter (fixed)"
$ws.Range("D5").Value  = "Tiger (fixed)"
$ws.Range("D6").Value  = "Tastar"

$ws.Range("D10").Value = "Future Ideas"
$ws.Range("D10").Font.Bold = $true

$ws.Range("D11").Value = "Graphical User Interface"
$ws.Range("D12").Value = "Link with PowerBI dashboard (automate)"
$ws.Range("D13").Value = "Scrub data over mastersheet to have uniform names/skus"
$ws.Range("D14").Value = "Create database"
$ws.Range("D15").Value = "Secondary sheet to simplify model"

$ws.Range("D21").Value = "Websites with Little Data"
$ws.Range("D21").Font.Bold = $true

$ws.Range("D22").Value = "Image Share"
$ws.Range("D22").Font.Bold = $false
$ws.Range("D23").Value = "Image Pro"
$ws.Range("D24").Value = "IP Store"
$ws.Range("D25").ClearContents()

# "Future Ideas" / "Websites with Little Data" headers no longer live in
# column A (row 22) -- unbold it now that it just reads "Zones".
$ws.Range("A22").Font.Bold = $false

# --- Column A: rest of the ProFocus websites --------------------------------
$ws.Range("A19").Value = "MacMall"
$ws.Range("A20").Value = "shi"
$ws.Range("A21").Value = "Grand & Toy"
$ws.Range("A22").Value = "Zones"
$ws.Range("A23").Value = "CDW"
$ws.Range("A24").Value = "IT Supplies"
$ws.Range("A25").Value = "Imaging Spectrum "
$ws.Range("A26").Value = "Laube"
$ws.Range("A27").Value = "LexJet"
$ws.Range("A28").Value = "Buffalo"
$ws.Range("A29").Value = "All American"
$ws.Range("A30").Value = "ProImaging Supplies"
$ws.Range("A31").Value = "Shades of Paper"
$ws.Range("A32").Value = "Spectraflow"

# --- Column D width grew slightly to fit the longer condensed strings ------
$ws.Columns.Item(4).ColumnWidth = 50.1666666666

# --- Selection moved to D6 before last save ---------------------------------
$ws.Range("D6").Select()
